{"js": "// 1. Title: \"Digital To Analog Converter\" -> \"Analog to Digital Converter\"\n//    (the diff shows the original single run split into \"Analog to Digital\" + \" Converter\";\n//     textually this is equivalent to replacing \"Digital To Analog\" with \"Analog to Digital\"\n//     and leaving \" Converter\" in place).\n{\n  const results = context.document.body.search(\"Digital To Analog\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Analog to Digital\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 2. Description paragraph: \"convertidor Digital a Analogico\" -> \"convertidor Analogico a Digital\"\n{\n  const results = context.document.body.search(\"Digital a Analogico\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Analogico a Digital\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 3. Wikipedia hyperlink: point to the English Analog-to-digital converter article\n//    (both the link target and the visible URL text change).\n{\n  const results = context.document.body.search(\n    \"https://es.wikipedia.org/wiki/Conversi%C3%B3n_digital-anal%C3%B3gica\",\n    { matchCase: true }\n  );\n  results.load(\"text,hyperlink\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const link = results.items[0];\n    link.hyperlink = \"https://en.wikipedia.org/wiki/Analog-to-digital_converter\";\n    await context.sync();\n    link.load(\"text,hyperlink\");\n    await context.sync();\n    link.insertText(\"https://en.wikipedia.org/wiki/Analog-to-digital_converter\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 4. Append a new \"Codigo de GitHub\" section with a link to the repo archive.\n{\n  const body = context.document.body;\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n\n  const lastParagraph = paras.items[paras.items.length - 1];\n  const githubLabelParagraph = lastParagraph.insertParagraph(\"Codigo de GitHub\", \"After\");\n  await context.sync();\n\n  const githubLinkParagraph = githubLabelParagraph.insertParagraph(\n    \"https://github.com/omaresl/SWEmbeddedTutorials/archive/ADC_Module_Example.zip\",\n    \"After\"\n  );\n  await context.sync();\n\n  const linkRange = githubLinkParagraph.getRange();\n  linkRange.hyperlink =\n    \"https://github.com/omaresl/SWEmbeddedTutorials/archive/ADC_Module_Example.zip\";\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Title: \"Digital To Analog Converter\" -> \"Analog to Digital Converter\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Digital To Analog Converter\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Analog to Digital Converter\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2. Description paragraph: \"convertidor Digital a Analogico\" -> \"convertidor Analogico a Digital\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Digital a Analogico\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Analogico a Digital\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# 3. Wikipedia hyperlink: point to the English Analog-to-digital converter article\n#    (update both the stored address and the visible display text).\nforeach ($h in $d.Hyperlinks) {\n    if ($h.Address -eq \"https://es.wikipedia.org/wiki/Conversi%C3%B3n_digital-anal%C3%B3gica\") {\n        $h.Address = \"https://en.wikipedia.org/wiki/Analog-to-digital_converter\"\n        $h.TextToDisplay = \"https://en.wikipedia.org/wiki/Analog-to-digital_converter\"\n    }\n}\n\n# 4. Append a new \"Codigo de GitHub\" section with a link to the repo archive.\n$endRange1 = $d.Content\n$endRange1.Collapse(0)   # wdCollapseEnd\n$endRange1.InsertParagraphAfter()\n$githubLabelPara = $d.Paragraphs($d.Paragraphs.Count)\n$githubLabelPara.Range.Text = \"Codigo de GitHub\"\n\n$endRange2 = $d.Content\n$endRange2.Collapse(0)\n$endRange2.InsertParagraphAfter()\n$githubLinkPara = $d.Paragraphs($d.Paragraphs.Count)\n$githubUrl = \"https://github.com/omaresl/SWEmbeddedTutorials/archive/ADC_Module_Example.zip\"\n$d.Hyperlinks.Add($githubLinkPara.Range, $githubUrl, [Type]::Missing, [Type]::Missing, $githubUrl) | Out-Null\n"}
